$d = $word.ActiveDocument

# --- Edit 1: "Implement organizations app models and services" --------
# Split the trailing word "services" into its own run and mark it red,
# leaving the rest of the sentence (still italic + green highlight) intact.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Implement*organizations*app*models*and*services*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $found = $r.Find.Execute("services", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Font.Color = 255
    }
}

# --- Edit 2: drop the trailing "Would you like me to elaborate..." ----
$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
if ($last.Range.Text -like "*Would you like me to elaborate*") {
    $last.Range.Delete()
}
